# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Pepino ensalada" (Comercializadora
# del Agro de Limari) right before the existing row 34, shifting the
# historical data down by two rows (and growing the used range from
# A1:R120 to A1:R122).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all the existing data rows (34..120) down by two rows.
$ws.Rows.Item(34).Insert()
$ws.Rows.Item(34).Insert()

# New row 34 - "Primera" quality, week of 2022-01-20
$ws.Range("A34").Value = 2
$ws.Range("B34").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").Value = 44581
$ws.Range("E34").Value = 4
$ws.Range("F34").Value = 100112043
$ws.Range("G34").Value = "Pepino ensalada"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 400
$ws.Range("K34").Value = 10500
$ws.Range("L34").Value = 11000
$ws.Range("M34").Value = 10750
$ws.Range("N34").Value = "$/caja 70 unidades"
$ws.Range("O34").Value = "Provincia de Limarí"
$ws.Range("P34").Value = 154
$ws.Range("Q34").Value = 70
$ws.Range("R34").Value = "Hortaliza"

# New row 35 - "Segunda" quality, week of 2022-01-20
$ws.Range("A35").Value = 2
$ws.Range("B35").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C35").Value = "Coquimbo"
$ws.Range("D35").Value = 44581
$ws.Range("E35").Value = 4
$ws.Range("F35").Value = 100112043
$ws.Range("G35").Value = "Pepino ensalada"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Segunda"
$ws.Range("J35").Value = 300
$ws.Range("K35").Value = 8500
$ws.Range("L35").Value = 9000
$ws.Range("M35").Value = 8750
$ws.Range("N35").Value = "$/caja 100 unidades"
$ws.Range("O35").Value = "Provincia de Limarí"
$ws.Range("P35").Value = 88
$ws.Range("Q35").Value = 100
$ws.Range("R35").Value = "Hortaliza"
